$wb = $excel.ActiveWorkbook

# Sheet "NBR" - update Reaction_number column (C) for rows 2-20
$wsNBR = $wb.Worksheets.Item("NBR")
$nbrValues = @(855, 851, 847, 843, 832, 835, 835, 835, 834, 829, 824, 819, 827, 824, 820, 823, 823, 823, 821)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# Sheet "BAR" - update Reaction_number column (C) for rows 2-20
$wsBAR = $wb.Worksheets.Item("BAR")
$barValues = @(749, 749, 745, 746, 764, 757, 756, 756, 751, 750, 747, 743, 725, 726, 726, 721, 721, 721, 721)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
